$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content change -------------------------------------------------------
# Row 7 / column D used to hold the combined task "validações data e
# dinheiro". It is being split into two separate, more specific validation
# tasks:
#   - "validações data"    -> new cell F6 (styled like the neighbouring E6)
#   - "validacao dinheiro" -> new cell D8 (styled like the neighbouring E6)
# and the old combined D7 cell is cleared out.

# F6: copy formatting from E6 (style 8) then set its text.
$ws.Range("E6").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F6").Value = "validações data"

# D8: copy formatting from E6 (style 8) then set its text.
$ws.Range("E6").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = "validacao dinheiro"

# Done with the clipboard-based format copies.
$excel.CutCopyMode = $false

# Clear the old combined text + formatting from D7, restoring it to the
# sheet's default (unstyled) look.
$ws.Range("D7").Clear()
$ws.Range("D7").HorizontalAlignment = 1

# The row also gained an extra (still empty) cell at E7.
$ws.Range("E7").HorizontalAlignment = 1

# --- Cosmetic re-format of the header rows ---------------------------------
# Rows 1-3, columns B:G now carry explicit (default) formatting too.
$ws.Range("B1:G3").HorizontalAlignment = 1

# --- View state -------------------------------------------------------------
# The saved view now scrolls so row 3 is the first visible row.
$ws.Application.ActiveWindow.ScrollRow = 3
